# Weapon balance sheet update (r260 -> r268)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header notes ---
$ws.Range("A1").Value = "Additional damage multiplieris applied for sniper rifles (2x)"
$ws.Range("K1").Value = "Last updated 16.3. (r268)"

# --- Move speed column (D) rebalanced per weapon category ---
# Handguns (rows 4-10): 120/70% -> 150/120%
$ws.Range("D4:D10").Value = "150/120%"

# Shotguns + SMGs (rows 11-21): 100/50% -> 125/110%
$ws.Range("D11:D21").Value = "125/110%"

# Assault Rifles + Thumper (rows 22-30): 80/30% -> 100/70%
$ws.Range("D22:D30").Value = "100/70%"

# Category-boundary rows also had their top/bottom border formatting cleared
$ws.Range("D17").ClearFormats()
$ws.Range("D17").Value = "125/110%"
$ws.Range("D21").ClearFormats()
$ws.Range("D21").Value = "125/110%"
$ws.Range("D29").ClearFormats()
$ws.Range("D29").Value = "100/70%"

# --- SMGs rows gain Fire rate (F) / Recoil (H) data ---
$ws.Range("F17").Value = "'--"
$ws.Range("H17").Value = "'--"

$ws.Range("F18").Value = "-"
$ws.Range("H18").Value = "-"

$ws.Range("F19").Value = "'+++"
$ws.Range("H19").Value = "'+++"

$ws.Range("F20").Value = "'++++"
$ws.Range("H20").Value = "'++++"

$ws.Range("F21").Value = "'-"
$ws.Range("H21").Value = "'++"

# --- ACR recoil correction ---
$ws.Range("H28").Value = "+"

# --- M240 clip/recoil correction ---
$ws.Range("E40").Value = 19
$ws.Range("H40").Value = "'++"

# --- Column D width now explicit ---
$ws.Columns.Item(4).ColumnWidth = 9.7109375

# --- Restore selection / scroll position ---
$ws.Range("M14").Select()
